$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark Word leaves behind from the
#    previous editing session (it sits right after the PERSON(...) row,
#    just before the grey "vlt. LIEBLINGSCOCKTAIL -> COCKTAIL" comment).
# ---------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# ---------------------------------------------------------------------
# 2. Insert the two new relation rows for the reworked ER-model right
#    before the "COCKTAILKAUF(...)" row:
#       HATLIEBLINGSCOCKTAIL(PID->PERSON, LCOCKID->LIEBLINGSCOCKTAIL)
#       LIEBLINGSCOCKTAIL(LIEBLINGSCOCKTAILID, COCKTAILNAME)
#    New paragraphs are created by asking the COCKTAILKAUF paragraph to
#    insert a sibling before itself, which gives us a blank paragraph
#    that already carries the right Consolas/24pt formatting; we then
#    fill in the text and colour/underline the runs.
# ---------------------------------------------------------------------

function Get-CocktailkaufParagraph {
    $rng = $d.Content
    $rng.Find.Execute("COCKTAILKAUF", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
    return $rng.Paragraphs(1)
}

# -- paragraph 1: HATLIEBLINGSCOCKTAIL(...) -----------------------------
$cocktailkauf = Get-CocktailkaufParagraph
$cocktailkauf.Range.InsertParagraphBefore()

$cocktailkauf = Get-CocktailkaufParagraph
$p1 = $cocktailkauf.Previous()
$p1.Range.Text = "HATLIEBLINGSCOCKTAIL(PID->PERSON, LCOCKID->LIEBLINGSCOCKTAIL)"
$p1.Range.Font.Color = 192

$u1 = $d.Content
$u1.Find.Execute("PID->PERSON, LCOCKID->LIEBLINGSCOCKTAIL", $true, $false, $false, $false, `
                  $false, $true, 1, $false, "", 0) | Out-Null
$u1.Font.Underline = 1

# -- paragraph 2: LIEBLINGSCOCKTAIL(...) --------------------------------
$cocktailkauf = Get-CocktailkaufParagraph
$cocktailkauf.Range.InsertParagraphBefore()

$cocktailkauf = Get-CocktailkaufParagraph
$p2 = $cocktailkauf.Previous()
$p2.Range.Text = "LIEBLINGSCOCKTAIL(LIEBLINGSCOCKTAILID, COCKTAILNAME)"
$p2.Range.Font.Color = 192

$u2 = $d.Content
$u2.Find.Execute("LIEBLINGSCOCKTAILID", $true, $false, $false, $false, `
                  $false, $true, 1, $false, "", 0) | Out-Null
$u2.Font.Underline = 1
